$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Constants ---
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Copy formatting from "donor" cells (existing styles) onto cells that
#    need a style they don't currently carry, BEFORE those donor cells (or
#    the target cells) get their own content/style changed later on.
# ---------------------------------------------------------------------------

# H14 needs the plain "right-thin + bottom-thin" box style that currently
# only lives on B6 - grab it first, before B6 gets cleared below.
$ws.Range("B6").Copy()
$ws.Range("H14").PasteSpecial($xlPasteFormats)

# A3 and G14 need the "left-thin + bottom-thin" style that lives on J3.
$ws.Range("J3").Copy()
$ws.Range("A3").PasteSpecial($xlPasteFormats)
$ws.Range("G14").PasteSpecial($xlPasteFormats)

# B3 needs the underlined "right-thin + bottom-thin" style that lives on K3.
$ws.Range("K3").Copy()
$ws.Range("B3").PasteSpecial($xlPasteFormats)

# G11 / H11 need the centered header-box styles that live on A1 / B1.
$ws.Range("A1").Copy()
$ws.Range("G11").PasteSpecial($xlPasteFormats)
$ws.Range("B1").Copy()
$ws.Range("H11").PasteSpecial($xlPasteFormats)

# G12 / H12 need the plain "Int/ID" row styles that live on A2 / B2.
$ws.Range("A2").Copy()
$ws.Range("G12").PasteSpecial($xlPasteFormats)
$ws.Range("B2").Copy()
$ws.Range("H12").PasteSpecial($xlPasteFormats)

# G13 / H13 need the plain "char/Name" row styles; use D3 / H3 as donors
# (A3/B3 are themselves being changed above).
$ws.Range("D3").Copy()
$ws.Range("G13").PasteSpecial($xlPasteFormats)
$ws.Range("H3").Copy()
$ws.Range("H13").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Update text content across the mockup tables.
# ---------------------------------------------------------------------------

# Table 1 header: "Account" -> "Kontakte"
$ws.Range("A1").Value = "Kontakte"

# Row 2 of table 1: ID -> User1, add "> Account.ID" in C2
$ws.Range("B2").Value = "User1"
$ws.Range("C2").Value = "> Account.ID"

# Row 3 of table 1: char/Name -> Int/User2, add "> Account.ID" in C3
$ws.Range("A3").Value = "Int"
$ws.Range("B3").Value = "User2"
$ws.Range("C3").Value = "> Account.ID"

# F3 reference stays "> Medien.ID" (unchanged text, only its shared-string
# index shifts, which the engine manages automatically).

# F4 reference relabelled from "> Account.ID" (index shift only, text same)

# Row 5 foreign-key label text stays "> Account.ID" (index shift only).

# New "Account" table placed at G11:H14
$ws.Range("G11").Value = "Account"
$ws.Range("G12").Value = "Int"
$ws.Range("H12").Value = "ID"
$ws.Range("G13").Value = "char"
$ws.Range("H13").Value = "Name"
$ws.Range("G14").Value = "char"
$ws.Range("H14").Value = "Passwort"

# ---------------------------------------------------------------------------
# 3) Remove obsolete rows / cells.
# ---------------------------------------------------------------------------

# Old "Passwort" row of table 1 is gone entirely.
$ws.Range("A4:B4").Clear()

# Old "Medien"/"Kontakte" rows of table 1 collapse into two blank separator
# rows (still boxed, no text, no foreign-key column).
$ws.Range("C5").Clear()
$ws.Range("C6").Clear()

# Footnotes under the old table 1 ("Medien und Kontakte...", and the
# "1. Medium erstellen..." reminder) are removed.
$ws.Range("A12").Clear()
$ws.Range("A13").Clear()

# ---------------------------------------------------------------------------
# 4) Turn A5/A6/B5/B6 into blank-but-formatted placeholder cells (no border,
#    explicit empty format so Excel keeps an <c s="..."/> entry for them).
# ---------------------------------------------------------------------------
foreach ($addr in @("A5", "B5", "A6", "B6")) {
    $r = $ws.Range($addr)
    $r.Value = ""
    $r.Borders.Item(7).LineStyle = -4142
    $r.Borders.Item(8).LineStyle = -4142
    $r.Borders.Item(9).LineStyle = -4142
    $r.Borders.Item(10).LineStyle = -4142
    $r.NumberFormat = "General"
}

# ---------------------------------------------------------------------------
# 5) Merge the new "Account" table header and size the new separator row.
# ---------------------------------------------------------------------------
$ws.Range("G11:H11").Merge()

$ws.Rows.Item(10).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 6) Selection as saved in the file.
# ---------------------------------------------------------------------------
$ws.Range("C4").Select()
